$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Random_Forest
$ws.Range("C3").Value = 94.56999999999999
$ws.Range("D3").Value = 93.52
$ws.Range("F3").Value = 95.28
$ws.Range("G3").Value = 94

# Row 4: Deep_Neural_Network
$ws.Range("C4").Value = 98.8
$ws.Range("D4").Value = 89.73999999999999
$ws.Range("E4").Value = 87.23
$ws.Range("F4").Value = 99.06
$ws.Range("G4").Value = 93.5
